$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fill in the last (empty) paragraph of the "8. Conclusiones"
#    section (the one immediately before the final sectPr) with the
#    author's conclusion text, reproducing the original run / proofErr
#    layout exactly.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">Mi conclusión acerca de las fuentes de información disponibles es que los documentos oficiales y contrastados están bastante limitados, y que tienes que recurrir a artículos de gente especializada en la materia para encontrar algo que sea lo suficientemente profesional. En la parte de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>OpenStreet</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Map</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> sí que es cierto que hay más documentos, pero como es una tecnología libre es más difícil de filtrar los artículos, mientras que en Google </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Maps</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> en la mayoría de casos tienes que recurrir a sus propias fuentes de información, lo que sesga una posible visión objetiva de la misma.</w:t>
  </w:r>
</w:p>
'@

[void]$r.InsertXML($xml)

# ------------------------------------------------------------------
# 2. The document grew by two pages once the conclusion text was
#    added, so the cached page-number field in the footer (previously
#    "6") now needs to read "8".
# ------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
foreach ($fld in $footer.Range.Fields) {
  if ($fld.Code.Text -match "PAGE") {
    $fld.Result.Text = "8"
  }
}
